# Insert a new data row at row 134 (pushing existing rows 134..169 down to
# 135..170) and populate it with the new price-report entry.
#
# Columns that are constant across every row of this sheet (Mercado ID,
# Mercado, Region, Codreg, Categoria ID, Categoria, Kg o Unidades,
# Clasificacion) are copied from the neighbouring row 133.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 134:169 down to 135:170, leaving a blank row 134 behind.
$ws.Rows(134).Insert()

$ws.Range("A134").Value = 1
$ws.Range("B134").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C134").Value = "Arica y Parinacota"
$ws.Range("D134").Value = 44985
$ws.Range("E134").Value = 15
$ws.Range("F134").Value = 100114001
$ws.Range("G134").Value = "Papa"
$ws.Range("H134").Value = "Cardinal"
$ws.Range("I134").Value = "1a (cosecha)"
$ws.Range("J134").Value = 950
$ws.Range("K134").Value = 14000
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = 14368
$ws.Range("N134").Value = "$/saco 25 kilos"
$ws.Range("O134").Value = "Región de Coquimbo"
$ws.Range("P134").Value = 575
$ws.Range("Q134").Value = 25
$ws.Range("R134").Value = "Hortaliza"
